# Product features.xlsx — apply the "image uploading" backlog update.
#
# 1. "product backlog" (sheet1), row 8: the organizational-chart logging
#    story is finished -> give it the same highlighted formatting as the
#    other closed rows and mark Status = "(3) Completed".
# 2. "sprint backlog" (sheet2): the old placeholder row (ID 8) is replaced
#    by the next sprint's stories - irs number field, user photos, org
#    logos - continuing the "image uploading" work mentioned in the commit.
# 3. View state: "sprint backlog" becomes the active/selected sheet, and
#    both sheets' selections move to where the author was last working.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # product backlog
$ws2 = $wb.Worksheets.Item(2)   # sprint backlog

# --- 1. product backlog: close out row 8 -----------------------------
$ws1.Rows.Item(8).Interior.Color = 65535   # same yellow fill as the other "done" rows
$ws1.Range("I8").Value = "(3) Completed"

# --- 2. sprint backlog: roll the plan forward -------------------------
$ws2.Range("A2").Value = 15
$ws2.Range("E2").Value = "Add irs number field to organization"
$ws2.Rows.Item(2).AutoFit()   # drop the old wrapped-text row height

$ws2.Range("A3").Value = 17
$ws2.Range("E3").Value = "add photo for users"

$ws2.Range("A4").Value = 18
$ws2.Range("E4").Value = "add logo for organizations"

# --- 3. selections / active sheet -------------------------------------
$ws1.Activate()
$ws1.Range("A17:XFD18").Select()

$ws2.Activate()
$ws2.Range("F9").Select()
